$d = $word.ActiveDocument

# Remove the three paragraphs that follow the "Project Plan" heading:
#   - an empty paragraph
#   - "Question:"
#   - "Analysis on the rise of expenditure on gym membership and supplements
#      to aid a portable personal trainer application."
# Paragraph index 2 is the first to go; paragraph index 4 is the last to go.
# Extend through the start of paragraph 5 so the paragraph marks themselves
# (and not just the run text) are removed, merging "Project Plan" directly
# into the formerly-5th paragraph.
$start = $d.Paragraphs.Item(2).Range.Start
$end = $d.Paragraphs.Item(5).Range.Start
$r = $d.Range($start, $end)
$r.Delete()

Write-Output "done"
